$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 17, shifting the
# existing rows 17-57 down to 19-59 (Excel carries formatting down with
# the shift automatically).
$ws.Range("A17:R18").Insert()

# Populate the two newly inserted rows with the new weekly record
# (market/category columns are constant across the sheet; only the
# date/quality/volume/price columns change).

# Row 17: Primera
$ws.Range("A17").Value = 11
$ws.Range("B17").Value = "Vega Monumental Concepción"
$ws.Range("C17").Value = "Bíobío"
$ws.Range("D17").Value = 44797
$ws.Range("E17").Value = 8
$ws.Range("F17").Value = 100112043
$ws.Range("G17").Value = "Pepino dulce"
$ws.Range("H17").Value = "Cultivar IV Región"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 50
$ws.Range("K17").Value = 15000
$ws.Range("L17").Value = 15000
$ws.Range("M17").Value = 15000
$ws.Range("N17").Value = "`$/bandeja 18 kilos"
$ws.Range("O17").Value = "Provincia de Limarí"
$ws.Range("P17").Value = 833
$ws.Range("Q17").Value = 18
$ws.Range("R17").Value = "Hortaliza"

# Row 18: Segunda
$ws.Range("A18").Value = 11
$ws.Range("B18").Value = "Vega Monumental Concepción"
$ws.Range("C18").Value = "Bíobío"
$ws.Range("D18").Value = 44797
$ws.Range("E18").Value = 8
$ws.Range("F18").Value = 100112043
$ws.Range("G18").Value = "Pepino dulce"
$ws.Range("H18").Value = "Cultivar IV Región"
$ws.Range("I18").Value = "Segunda"
$ws.Range("J18").Value = 50
$ws.Range("K18").Value = 12000
$ws.Range("L18").Value = 12000
$ws.Range("M18").Value = 12000
$ws.Range("N18").Value = "`$/bandeja 18 kilos"
$ws.Range("O18").Value = "Provincia de Limarí"
$ws.Range("P18").Value = 667
$ws.Range("Q18").Value = 18
$ws.Range("R18").Value = "Hortaliza"
